$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9 through 31 (old trailing rows no longer needed)
$ws.Range("A9:A31").EntireRow.Delete() | Out-Null

# Update A2:A8 with the new combined tuple-like strings
$ws.Range("A2").Value = "('Giant Growth', ['{G}', 'Instant', 'Target creature gets +3/+3 until end of turn.'])"
$ws.Range("A3").Value = "('Hinder', ['{1}{U}{U}', 'Instant', 'Counter target spell. If that spell is countered this way, put that card on the top or bottom of its owner" + [char]0x2019 + "s library instead of into that player" + [char]0x2019 + "s graveyard.'])"
$ws.Range("A4").Value = "('Hypnotic Specter', ['{1}{B}{B}', 'Creature " + [char]0x2014 + " Specter', 'Flying', 'Whenever Hypnotic Specter deals damage to an opponent, that player discards a card at random.', '2/2'])"
$ws.Range("A5").Value = "('Lightning Helix', ['{R}{W}', 'Instant', 'Lightning Helix deals 3 damage to any target and you gain 3 life.'])"
$ws.Range("A6").Value = "('Putrefy', ['{1}{B}{G}', 'Instant', 'Destroy target artifact or creature. It can" + [char]0x2019 + "t be regenerated.'])"
$ws.Range("A7").Value = "('Pyroclasm', ['{1}{R}', 'Sorcery', 'Pyroclasm deals 2 damage to each creature.'])"
$ws.Range("A8").Value = "('Zombify', ['{3}{B}', 'Sorcery', 'Return target creature card from your graveyard to the battlefield.'])"
